# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de
# handback has completed (zh-cn was already handed back earlier):
#   - "Ready for handoff" -> "Handed back: in sync with en-US" on the
#     Overview sheet status columns
#   - Populates "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" on both the zh-cn and de-de detail sheets
#   - Widens the now-longer text columns to fit the new content

$wb = $excel.ActiveWorkbook

$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e389bb62f6334bf2afa44b0a57ca0d4a0330c03/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: status cells + column widths
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-15 16:35:23"
$zhcn.Range("K3").Value = "2016-08-15 16:35:23"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aUrl, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aUrl, "", "", "a.md")

$zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("K2").Value = "2016-08-15 16:35:30"
$dede.Range("K3").Value = "2016-08-15 16:35:30"

$dede.Hyperlinks.Add($dede.Range("I2"), $aUrl, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $aUrl, "", "", "a.md")

$dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$dede.Columns.Item(10).ColumnWidth = 39.1666666666667
